# Updated cryptos list on Fri Sep 27 09:58:21 UTC 2024 with GitHub Actions
#
# All Price/Coin/Link cells in this sheet are stored as plain text (General
# number format), even when the text looks like a number (e.g. "1.00",
# "5.90"). Assigning such a string straight to Range.Value lets Excel's
# automatic type inference silently coerce it to a numeric cell, which would
# diverge from the source data. Set-TextValue guards against that: it flips
# the cell to text format, writes the literal string, then restores the
# "Normal" cell style so the cell ends up back on the sheet's default
# (unstyled/General) formatting - matching cells that never needed the
# text-format workaround (e.g. multi-dot price strings, which Excel can't
# parse as numbers anyway and so are never coerced).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "65.802.94"
Set-TextValue $ws.Range("E2") "  +3.17%  "
Set-TextValue $ws.Range("D3") "2.666.84"
Set-TextValue $ws.Range("E3") "  +1.82%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "606.77"
Set-TextValue $ws.Range("E5") "  +2.08%  "
Set-TextValue $ws.Range("D6") "158.13"
Set-TextValue $ws.Range("E6") "  +4.84%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("D8") "0.589"
Set-TextValue $ws.Range("E8") "  +0.36%  "
Set-TextValue $ws.Range("D9") "2.664.94"
Set-TextValue $ws.Range("E9") "  +1.79%  "
Set-TextValue $ws.Range("E10") "  +9.35%  "
Set-TextValue $ws.Range("D11") "0.404"
Set-TextValue $ws.Range("E11") "  +2.72%  "
Set-TextValue $ws.Range("D12") "5.90"
Set-TextValue $ws.Range("E12") "  +1.81%  "
Set-TextValue $ws.Range("E13") "  +1.75%  "
Set-TextValue $ws.Range("D14") "29.99"
Set-TextValue $ws.Range("E14") "  +7.72%  "
Set-TextValue $ws.Range("D15") "0.0000197"
Set-TextValue $ws.Range("E15") "  +16.29%  "
Set-TextValue $ws.Range("D16") "3.147.13"
Set-TextValue $ws.Range("E16") "  +1.80%  "
Set-TextValue $ws.Range("D17") "65.538.31"
Set-TextValue $ws.Range("E17") "  +2.88%  "
Set-TextValue $ws.Range("D18") "2.649.79"
Set-TextValue $ws.Range("E18") "  +0.62%  "
Set-TextValue $ws.Range("D19") "12.72"
Set-TextValue $ws.Range("E19") "  +4.04%  "
Set-TextValue $ws.Range("D20") "4.90"
Set-TextValue $ws.Range("E20") "  +2.72%  "
Set-TextValue $ws.Range("D21") "360.98"
Set-TextValue $ws.Range("E21") "  +3.79%  "
Set-TextValue $ws.Range("D22") "7.47"
Set-TextValue $ws.Range("E22") "  +6.70%  "
Set-TextValue $ws.Range("E23") "  +0.10%  "
Set-TextValue $ws.Range("D24") "69.42"
Set-TextValue $ws.Range("E24") "  +2.99%  "
Set-TextValue $ws.Range("D25") "1.70"
Set-TextValue $ws.Range("E25") "  +0.96%  "
Set-TextValue $ws.Range("D26") "9.53"
Set-TextValue $ws.Range("E26") "  +4.08%  "
Set-TextValue $ws.Range("E27") "  +17.66%  "
Set-TextValue $ws.Range("D28") "1.64"
Set-TextValue $ws.Range("E28") "  -1.96%  "
Set-TextValue $ws.Range("E29") "  +2.78%  "
Set-TextValue $ws.Range("D30") "8.19"
Set-TextValue $ws.Range("E30") "  -0.99%  "
Set-TextValue $ws.Range("E31") "  +6.32%  "
Set-TextValue $ws.Range("B32") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D32") "1.00"
Set-TextValue $ws.Range("E32") "  +0.16%  "
Set-TextValue $ws.Range("B33") "Bittensor"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D33") "539.53"
Set-TextValue $ws.Range("E33") "  -1.35%  "
Set-TextValue $ws.Range("D34") "1.81"
Set-TextValue $ws.Range("E34") "  -0.55%  "
Set-TextValue $ws.Range("E35") "  +2.86%  "
Set-TextValue $ws.Range("D36") "6.37"
Set-TextValue $ws.Range("E36") "  +3.80%  "
Set-TextValue $ws.Range("E37") "  +3.30%  "
Set-TextValue $ws.Range("D38") "20.78"
Set-TextValue $ws.Range("E38") "  +4.34%  "
Set-TextValue $ws.Range("E39") "  +1.94%  "
Set-TextValue $ws.Range("D40") "162.63"
Set-TextValue $ws.Range("E40") "  -1.40%  "
Set-TextValue $ws.Range("E41") "  +0.01%  "
Set-TextValue $ws.Range("E42") "  +0.03%  "
Set-TextValue $ws.Range("D43") "42.46"
Set-TextValue $ws.Range("E43") "  +6.25%  "
Set-TextValue $ws.Range("D44") "166.41"
Set-TextValue $ws.Range("E44") "  -1.08%  "
Set-TextValue $ws.Range("E45") "  +2.07%  "
Set-TextValue $ws.Range("E46") "  +9.39%  "
Set-TextValue $ws.Range("E47") "  +5.45%  "
Set-TextValue $ws.Range("D48") "23.09"
Set-TextValue $ws.Range("E48") "  -0.30%  "
Set-TextValue $ws.Range("D49") "0.661"
Set-TextValue $ws.Range("E49") "  +3.94%  "
Set-TextValue $ws.Range("E50") "  +5.35%  "
Set-TextValue $ws.Range("E51") "  +1.80%  "
